{"js": "// Update the \"compatible with\" blurb to reference the Fall Creators Update\n// SDK instead of the Creators Update SDK, matching the commit:\n// \"Updated for Windows 10 Fall Creators Update SDK (16299)\".\n\nconst body = context.document.body;\nconst oldText =\n  \"This sample is compatible with the Windows 10 Creators Update SDK (15063)\";\nconst newText =\n  \"This sample is compatible with the Windows 10 Fall Creators Update SDK (16299)\";\n\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the SDK compatibility sentence.\");\n}\n\nconst targetRange = results.items[0];\nconst paras = targetRange.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nconst targetPara = paras.items[0];\n\n// Rewrite the sentence in place; the surrounding run keeps its italic\n// formatting since we're replacing text inside the existing range.\ntargetRange.insertText(newText, Word.InsertLocation.replace);\n\n// The paragraph that used to follow only held the invisible \"_GoBack\"\n// bookmark (no visible text). Word folds its paragraph mark into the\n// sentence above and opens a fresh blank paragraph afterwards -- reproduce\n// that by dropping the next paragraph and inserting a new empty one.\nconst nextPara = targetPara.getNextOrNullObject();\nawait context.sync();\nif (!nextPara.isNullObject) {\n  nextPara.delete();\n  targetPara.insertParagraph(\"\", Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Update the \"compatible with\" blurb to reference the Fall Creators Update\n# SDK instead of the Creators Update SDK, matching the commit:\n# \"Updated for Windows 10 Fall Creators Update SDK (16299)\".\n\n$d = $word.ActiveDocument\n\n$oldText = \"This sample is compatible with the Windows 10 Creators Update SDK (15063)\"\n$newText = \"This sample is compatible with the Windows 10 Fall Creators Update SDK (16299)\"\n\n$find = $d.Content.Find\n$find.Text = $oldText\n$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n\n# Re-locate the (now updated) paragraph's 1-based index.\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $newText) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -gt 0 -and $targetIndex -lt $d.Paragraphs.Count) {\n    # The paragraph right after used to contain only the invisible\n    # \"_GoBack\" bookmark. Word folds its paragraph mark into the sentence\n    # above and opens a fresh blank paragraph afterwards -- reproduce that\n    # by deleting the next paragraph and inserting a new empty one.\n    $nextPara = $d.Paragraphs.Item($targetIndex + 1)\n    $nextPara.Range.Delete()\n\n    $targetPara = $d.Paragraphs.Item($targetIndex)\n    $targetPara.Range.InsertParagraphAfter()\n}\n"}
